$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 352 (Fukumoto/Primera),
# pushing the existing rows 352-359 down to 354-361.
$ws.Rows.Item(352).Insert()
$ws.Rows.Item(352).Insert()

# New row 352: Naranja - Valencia - Primera
$ws.Cells.Item(352, 1).Value = 7
$ws.Cells.Item(352, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(352, 3).Value = "Ñuble"
$ws.Cells.Item(352, 4).Value = 44656
$ws.Cells.Item(352, 5).Value = 16
$ws.Cells.Item(352, 6).Value = "Fruta"
$ws.Cells.Item(352, 7).Value = 100102
$ws.Cells.Item(352, 8).Value = "Cítricos"
$ws.Cells.Item(352, 9).Value = 100102005
$ws.Cells.Item(352, 10).Value = "Naranja"
$ws.Cells.Item(352, 11).Value = "Valencia"
$ws.Cells.Item(352, 12).Value = "Primera"
$ws.Cells.Item(352, 13).Value = 120
$ws.Cells.Item(352, 14).Value = 9500
$ws.Cells.Item(352, 15).Value = 10000
$ws.Cells.Item(352, 16).Value = 9750
$ws.Cells.Item(352, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(352, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(352, 19).Value = 650
$ws.Cells.Item(352, 20).Value = 15

# New row 353: Naranja - Valencia - Segunda
$ws.Cells.Item(353, 1).Value = 7
$ws.Cells.Item(353, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(353, 3).Value = "Ñuble"
$ws.Cells.Item(353, 4).Value = 44656
$ws.Cells.Item(353, 5).Value = 16
$ws.Cells.Item(353, 6).Value = "Fruta"
$ws.Cells.Item(353, 7).Value = 100102
$ws.Cells.Item(353, 8).Value = "Cítricos"
$ws.Cells.Item(353, 9).Value = 100102005
$ws.Cells.Item(353, 10).Value = "Naranja"
$ws.Cells.Item(353, 11).Value = "Valencia"
$ws.Cells.Item(353, 12).Value = "Segunda"
$ws.Cells.Item(353, 13).Value = 100
$ws.Cells.Item(353, 14).Value = 8000
$ws.Cells.Item(353, 15).Value = 8500
$ws.Cells.Item(353, 16).Value = 8250
$ws.Cells.Item(353, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(353, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(353, 19).Value = 550
$ws.Cells.Item(353, 20).Value = 15
